$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 154603
$ws.Range("C4").Value = 145733
$ws.Range("C5").Value = 8871
$ws.Range("C8").Value = 63.5
